$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D text-like numeric values remain exact text (not auto-converted to numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.188.46"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.685.89"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.04"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.18"
$ws.Range("E8").Value = "  +7.99%  "
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0890"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.924.78"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.692.47"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.93"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.191.83"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.02"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.61"
$ws.Range("E23").Value = "  +3.75%  "
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.57"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.33"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.46"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.541.59"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  +2.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.946"
$ws.Range("E37").Value = "  +2.65%  "
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.15"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.832.35"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.793"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.16"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").Value = "  +5.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.61"
$ws.Range("E49").Value = "  +5.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.29"
$ws.Range("E50").Value = "  +5.40%  "
$ws.Range("E51").Value = "  -1.80%  "
